# Auto-generated profit-column refresh, mirrors a scheduled market-data pull.
$wb = $excel.ActiveWorkbook

# ---- ALC: 23 cell updates ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 41.75
$ws.Range("I5").Value = 50.666668
$ws.Range("K5").Value = 50.666668
$ws.Range("M5").Value = 64.333332
$ws.Range("H19").Value = 1124.5
$ws.Range("J19").Value = 792.1429000000001
$ws.Range("L19").Value = 792.1429000000001
$ws.Range("N19").Value = -1142.1429
$ws.Range("H100").Value = 4999.4243
$ws.Range("I100").Value = 3655.85
$ws.Range("J100").Value = 7066.4614
$ws.Range("K100").Value = 3655.85
$ws.Range("L100").Value = 7066.4614
$ws.Range("M100").Value = -3114.85
$ws.Range("N100").Value = -8148.4614
$ws.Range("H121").Value = 2831.6667
$ws.Range("J121").Value = 2831.6667
$ws.Range("L121").Value = 8495.000100000001
$ws.Range("N121").Value = -11989.0001
$ws.Range("H141").Value = 5981.647
$ws.Range("I141").Value = 4445.8667
$ws.Range("K141").Value = 13337.6001
$ws.Range("M141").Value = -8157.6001

# ---- ARM: 43 cell updates ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4135.544
$ws.Range("I32").Value = 4590
$ws.Range("K32").Value = 4590
$ws.Range("M32").Value = -4303
$ws.Range("H45").Value = 2317
$ws.Range("I45").Value = 2022.52
$ws.Range("K45").Value = 2022.52
$ws.Range("M45").Value = -1645.52
$ws.Range("H61").Value = 5457
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 6066.3335
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 6066.3335
$ws.Range("M61").Value = -4788
$ws.Range("N61").Value = -6490.3335
$ws.Range("H74").Value = 1304.7693
$ws.Range("I74").Value = 1217.3226
$ws.Range("J74").Value = 1643.625
$ws.Range("K74").Value = 1217.3226
$ws.Range("L74").Value = 1643.625
$ws.Range("M74").Value = -343.3226
$ws.Range("N74").Value = -3391.625
$ws.Range("H77").Value = 1304.7693
$ws.Range("I77").Value = 1217.3226
$ws.Range("J77").Value = 1643.625
$ws.Range("K77").Value = 6086.612999999999
$ws.Range("L77").Value = 8218.125
$ws.Range("M77").Value = -1718.612999999999
$ws.Range("N77").Value = -16954.125
$ws.Range("H97").Value = 2244.6
$ws.Range("I97").Value = 1091.3572
$ws.Range("J97").Value = 4935.5
$ws.Range("K97").Value = 1091.3572
$ws.Range("L97").Value = 4935.5
$ws.Range("M97").Value = -595.3571999999999
$ws.Range("N97").Value = -5927.5
$ws.Range("H136").Value = 5457
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 6066.3335
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 18199.0005
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -23299.0005

# ---- BSM: 35 cell updates ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 47500
$ws.Range("J50").Value = 47500
$ws.Range("L50").Value = 47500
$ws.Range("N50").Value = -48648
$ws.Range("H94").Value = 2214.3333
$ws.Range("I94").Value = 836.46155
$ws.Range("K94").Value = 836.46155
$ws.Range("M94").Value = -385.46155
$ws.Range("H99").Value = 26855
$ws.Range("I99").Value = 1046.9231
$ws.Range("K99").Value = 1046.9231
$ws.Range("M99").Value = 451.0769
$ws.Range("H102").Value = 11537.625
$ws.Range("I102").Value = 11537.625
$ws.Range("K102").Value = 11537.625
$ws.Range("M102").Value = -8292.625
$ws.Range("H105").Value = 5174
$ws.Range("I105").Value = 1996
$ws.Range("J105").Value = 6233.3335
$ws.Range("K105").Value = 1996
$ws.Range("L105").Value = 6233.3335
$ws.Range("M105").Value = -249
$ws.Range("N105").Value = -9727.333500000001
$ws.Range("H122").Value = 79000
$ws.Range("J122").Value = 79000
$ws.Range("L122").Value = 79000
$ws.Range("N122").Value = -88800
$ws.Range("H134").Value = 2204.2
$ws.Range("I134").Value = 1255.25
$ws.Range("K134").Value = 3765.75
$ws.Range("M134").Value = -1230.75
$ws.Range("H138").Value = 58234.883
$ws.Range("J138").Value = 58234.883
$ws.Range("L138").Value = 58234.883
$ws.Range("N138").Value = -68514.883

# ---- CRP: 55 cell updates ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2377.818
$ws.Range("I31").Value = 1884.3334
$ws.Range("K31").Value = 1884.3334
$ws.Range("M31").Value = -1589.3334
$ws.Range("H34").Value = 2377.818
$ws.Range("I34").Value = 1884.3334
$ws.Range("K34").Value = 1884.3334
$ws.Range("M34").Value = -1682.3334
$ws.Range("H58").Value = 1616.8235
$ws.Range("I58").Value = 1461.5
$ws.Range("J58").Value = 1754.8889
$ws.Range("K58").Value = 1461.5
$ws.Range("L58").Value = 1754.8889
$ws.Range("M58").Value = -1258.5
$ws.Range("N58").Value = -2160.8889
$ws.Range("H70").Value = 35009.5
$ws.Range("I70").Value = 20081
$ws.Range("J70").Value = 37142.145
$ws.Range("K70").Value = 20081
$ws.Range("L70").Value = 37142.145
$ws.Range("M70").Value = -19766
$ws.Range("N70").Value = -37772.145
$ws.Range("H73").Value = 35009.5
$ws.Range("I73").Value = 20081
$ws.Range("J73").Value = 37142.145
$ws.Range("K73").Value = 20081
$ws.Range("L73").Value = 37142.145
$ws.Range("M73").Value = -18989
$ws.Range("N73").Value = -39326.145
$ws.Range("H94").Value = 1396.3334
$ws.Range("I94").Value = 983.8
$ws.Range("J94").Value = 1691
$ws.Range("K94").Value = 983.8
$ws.Range("L94").Value = 1691
$ws.Range("M94").Value = -532.8
$ws.Range("N94").Value = -2593
$ws.Range("H99").Value = 31912930
$ws.Range("I99").Value = 10368313
$ws.Range("K99").Value = 10368313
$ws.Range("M99").Value = -10366815
$ws.Range("H126").Value = 31912930
$ws.Range("I126").Value = 10368313
$ws.Range("K126").Value = 31104939
$ws.Range("M126").Value = -31102469
$ws.Range("H132").Value = 4855.4614
$ws.Range("I132").Value = 3715.7144
$ws.Range("K132").Value = 11147.1432
$ws.Range("M132").Value = -8617.143199999999
$ws.Range("H136").Value = 1616.8235
$ws.Range("I136").Value = 1461.5
$ws.Range("J136").Value = 1754.8889
$ws.Range("K136").Value = 4384.5
$ws.Range("L136").Value = 5264.6667
$ws.Range("M136").Value = -1834.5
$ws.Range("N136").Value = -10364.6667

# ---- CUL: 18 cell updates ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1715.875
$ws.Range("J34").Value = 2250.5
$ws.Range("L34").Value = 6751.5
$ws.Range("N34").Value = -6919.5
$ws.Range("H55").Value = 34335664
$ws.Range("I55").Value = 1599.6666
$ws.Range("J55").Value = 49050264
$ws.Range("K55").Value = 4798.9998
$ws.Range("L55").Value = 147150792
$ws.Range("M55").Value = -4621.9998
$ws.Range("N55").Value = -147151146
$ws.Range("H131").Value = 2189.842
$ws.Range("I131").Value = 1411.5555
$ws.Range("J131").Value = 2890.3
$ws.Range("K131").Value = 4234.666499999999
$ws.Range("L131").Value = 8670.900000000001
$ws.Range("M131").Value = 805.3335000000006
$ws.Range("N131").Value = -18750.9

# ---- GSM: 19 cell updates ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 109513.73
$ws.Range("I70").Value = 191184.83
$ws.Range("K70").Value = 191184.83
$ws.Range("M70").Value = -190914.83
$ws.Range("H73").Value = 109513.73
$ws.Range("I73").Value = 191184.83
$ws.Range("K73").Value = 191184.83
$ws.Range("M73").Value = -190248.83
$ws.Range("H122").Value = 1930.4584
$ws.Range("I122").Value = 1635.3572
$ws.Range("J122").Value = 2343.6
$ws.Range("K122").Value = 4906.071599999999
$ws.Range("L122").Value = 7030.799999999999
$ws.Range("M122").Value = -2456.071599999999
$ws.Range("N122").Value = -11930.8
$ws.Range("H132").Value = 6098.6816
$ws.Range("I132").Value = 5652.6665
$ws.Range("K132").Value = 16957.9995
$ws.Range("M132").Value = -14427.9995

# ---- LTW: 25 cell updates ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1062.5
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1166.6666
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 1166.6666
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -2664.6666
$ws.Range("H71").Value = 1062.5
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1166.6666
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 5833.333000000001
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -13321.333
$ws.Range("H93").Value = 2600
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 3200
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 3200
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -5696
$ws.Range("H109").Value = 21099.8
$ws.Range("J109").Value = 21099.8
$ws.Range("L109").Value = 21099.8
$ws.Range("N109").Value = -23873.8

# ---- WVR: 38 cell updates ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 74845.57000000001
$ws.Range("I96").Value = 102782.2
$ws.Range("K96").Value = 102782.2
$ws.Range("M96").Value = -101409.2
$ws.Range("H113").Value = 506
$ws.Range("I113").Value = 408.16666
$ws.Range("J113").Value = 799.5
$ws.Range("K113").Value = 1224.49998
$ws.Range("L113").Value = 2398.5
$ws.Range("M113").Value = 945.5000199999999
$ws.Range("N113").Value = -6738.5
$ws.Range("H122").Value = 2452.9285
$ws.Range("I122").Value = 2414.4
$ws.Range("J122").Value = 2549.25
$ws.Range("K122").Value = 7243.200000000001
$ws.Range("L122").Value = 7647.75
$ws.Range("M122").Value = -4793.200000000001
$ws.Range("N122").Value = -12547.75
$ws.Range("H135").Value = 50333.332
$ws.Range("J135").Value = 50333.332
$ws.Range("L135").Value = 50333.332
$ws.Range("N135").Value = -60473.332
$ws.Range("H136").Value = 2452.25
$ws.Range("I136").Value = 2452.25
$ws.Range("K136").Value = 7356.75
$ws.Range("M136").Value = -4806.75
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 62498
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
